# Re-applies the latest cryptos.xlsx price/volume snapshot (GitHub Actions run).
# Price cells that look numeric ("1.001", "0.3713", ...) are written with a
# leading apostrophe so Excel keeps them as literal text (matching the original
# inline-string cells) instead of coercing them into Doubles.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.463.35'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '1.574.06'
$ws.Range('E3').Value = '  +0.80%  '
$ws.Range('D4').Value = '''1.001'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''1.001'
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').Value = '''287.89'
$ws.Range('E6').Value = '  +0.66%  '
$ws.Range('D7').Value = '''0.3713'
$ws.Range('E7').Value = '  +1.91%  '
$ws.Range('D8').Value = '''47.58'
$ws.Range('E8').Value = '  -1.48%  '
$ws.Range('E9').Value = '  -0.29%  '
$ws.Range('D10').Value = '''1.149'
$ws.Range('E10').Value = '  +2.06%  '
$ws.Range('D11').Value = '''0.07557'
$ws.Range('E11').Value = '  +2.06%  '
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('D13').Value = '''20.79'
$ws.Range('E13').Value = '  -0.01%  '
$ws.Range('D14').Value = '''5.952'
$ws.Range('E14').Value = '  +0.59%  '
$ws.Range('D15').Value = '''6.944'
$ws.Range('E15').Value = '  +0.98%  '
$ws.Range('D16').Value = '1.568.50'
$ws.Range('E16').Value = '  +0.46%  '
$ws.Range('D17').Value = '''0.00001122'
$ws.Range('E17').Value = '  +1.79%  '
$ws.Range('D18').Value = '''88.22'
$ws.Range('E18').Value = '  -0.45%  '
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('D20').Value = '''6.407'
$ws.Range('E20').Value = '  +1.40%  '
$ws.Range('D21').Value = '''1.000'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').Value = '''16.54'
$ws.Range('E22').Value = '  +3.18%  '
$ws.Range('D23').Value = '''12.03'
$ws.Range('E23').Value = '  +0.85%  '
$ws.Range('D24').Value = '22.460.24'
$ws.Range('E24').Value = '  +0.27%  '
$ws.Range('D25').Value = '''2.391'
$ws.Range('E25').Value = '  -0.64%  '
$ws.Range('D26').Value = '''2.628'
$ws.Range('E26').Value = '  +3.59%  '
$ws.Range('D27').Value = '''151.34'
$ws.Range('E27').Value = '  +1.30%  '
$ws.Range('D28').Value = '''19.67'
$ws.Range('E28').Value = '  +1.57%  '
$ws.Range('D29').Value = '''4.964'
$ws.Range('E29').Value = '  -0.67%  '
$ws.Range('D30').Value = '''125.39'
$ws.Range('E30').Value = '  +2.16%  '
$ws.Range('D31').Value = '1.760.55'
$ws.Range('E31').Value = '  +1.44%  '
$ws.Range('D32').Value = '''1.093'
$ws.Range('E32').Value = '  +2.98%  '
$ws.Range('D33').Value = '''6.106'
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('D34').Value = '''1.991'
$ws.Range('E34').Value = '  +0.72%  '
$ws.Range('D35').Value = '''9.883'
$ws.Range('E35').Value = '  +3.10%  '
$ws.Range('D36').Value = '''0.08357'
$ws.Range('E36').Value = '  +1.55%  '
$ws.Range('D37').Value = '''0.02467'
$ws.Range('E37').Value = '  +3.40%  '
$ws.Range('D38').Value = '''0.2239'
$ws.Range('E38').Value = '  +1.56%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = '''1.307'
$ws.Range('E39').Value = '  +0.78%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = '''0.06385'
$ws.Range('E40').Value = '  +0.30%  '
$ws.Range('D41').Value = '''5.358'
$ws.Range('E41').Value = '  +0.51%  '
$ws.Range('D42').Value = '''11.46'
$ws.Range('E42').Value = '  +3.06%  '
$ws.Range('D43').Value = '''0.6275'
$ws.Range('E43').Value = '  +3.61%  '
$ws.Range('D44').Value = '''14.07'
$ws.Range('E44').Value = '  +3.52%  '
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').Value = '''0.6101'
$ws.Range('E46').Value = '  +6.41%  '
$ws.Range('D47').Value = '''3.777'
$ws.Range('E47').Value = '  +0.51%  '
$ws.Range('D48').Value = '''2.053'
$ws.Range('E48').Value = '  +2.56%  '
$ws.Range('D49').Value = '''125.14'
$ws.Range('E49').Value = '  +0.56%  '
$ws.Range('D50').Value = '''1.212'
$ws.Range('E50').Value = '  +0.39%  '
$ws.Range('D51').Value = '''0.07217'
$ws.Range('E51').Value = '  +0.01%  '
